$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 278189
$ws.Range("D10").Value = 29589
$ws.Range("E10").Value = 1752190865

$ws.Range("C13").Value = 37909
$ws.Range("E13").Value = 97626384

$ws.Range("C19").Value = 108917
$ws.Range("E19").Value = 344661124

$ws.Range("C63").Value = 40863
$ws.Range("E63").Value = 116446014

$ws.Range("C99").Value = 136567
$ws.Range("E99").Value = 863020929

$ws.Range("C103").Value = 48903
$ws.Range("E103").Value = 138395825

$ws.Range("C115").Value = 17547
$ws.Range("E115").Value = 38603446

$ws.Range("C117").Value = 19700
$ws.Range("E117").Value = 56409367

$ws.Range("C121").Value = 5962
$ws.Range("E121").Value = 11514574

$ws.Range("C122").Value = 9689
$ws.Range("E122").Value = 31925215

$ws.Range("C168").Value = 284917
$ws.Range("E168").Value = 1208113638

$ws.Range("C170").Value = 367253
$ws.Range("D170").Value = 38108
$ws.Range("E170").Value = 2843759227

$ws.Range("C171").Value = 115097
$ws.Range("D171").Value = 20263
$ws.Range("E171").Value = 444575953

$ws.Range("C173").Value = 54381
$ws.Range("E173").Value = 151843423

$ws.Range("C174").Value = 357149
$ws.Range("E174").Value = 1016456243

$ws.Range("C175").Value = 125501
$ws.Range("E175").Value = 811321192

$ws.Range("C179").Value = 235645
$ws.Range("E179").Value = 812018232

$ws.Range("C188").Value = 19700
$ws.Range("E188").Value = 66017381

$ws.Range("C203").Value = 13100
$ws.Range("E203").Value = 32991663

$ws.Range("C204").Value = 4750
$ws.Range("E204").Value = 11633084

$ws.Range("C205").Value = 11117
$ws.Range("E205").Value = 44051783

$ws.Range("C210").Value = 6418
$ws.Range("E210").Value = 18550633

$ws.Range("C213").Value = 3630
$ws.Range("E213").Value = 11065449

$ws.Range("C220").Value = 4712
$ws.Range("E220").Value = 11657490

$ws.Range("C262").Value = 38982
$ws.Range("E262").Value = 124708750

$ws.Range("C276").Value = 216630
$ws.Range("E276").Value = 1209898249
